$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 991.2857
$ws.Range("I98").Value = 992.3333
$ws.Range("K98").Value = 992.3333
$ws.Range("M98").Value = 505.6667
$ws.Range("H112").Value = 3216.074
$ws.Range("I112").Value = 997.5
$ws.Range("J112").Value = 3393.56
$ws.Range("K112").Value = 2992.5
$ws.Range("L112").Value = 10180.68
$ws.Range("M112").Value = -1884.5
$ws.Range("N112").Value = -12396.68
$ws.Range("H122").Value = 991.2857
$ws.Range("I122").Value = 992.3333
$ws.Range("K122").Value = 2976.9999
$ws.Range("M122").Value = -526.9998999999998
$ws.Range("H132").Value = 586.8276
$ws.Range("I132").Value = 566.7406999999999
$ws.Range("K132").Value = 1700.2221
$ws.Range("M132").Value = 829.7779
$ws.Range("H138").Value = 2861.652
$ws.Range("J138").Value = 3895.2856
$ws.Range("L138").Value = 11685.8568
$ws.Range("N138").Value = -21965.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1654.0588
$ws.Range("I32").Value = 1405.2449
$ws.Range("K32").Value = 1405.2449
$ws.Range("M32").Value = -1118.2449
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H61").Value = 2541.1765
$ws.Range("J61").Value = 4228.4287
$ws.Range("L61").Value = 4228.4287
$ws.Range("N61").Value = -4652.4287
$ws.Range("H136").Value = 2541.1765
$ws.Range("J136").Value = 4228.4287
$ws.Range("L136").Value = 12685.2861
$ws.Range("N136").Value = -17785.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2748.5
$ws.Range("I16").Value = 2631.5
$ws.Range("J16").Value = 3099.5
$ws.Range("K16").Value = 2631.5
$ws.Range("L16").Value = 3099.5
$ws.Range("M16").Value = -2344.5
$ws.Range("N16").Value = -3673.5
$ws.Range("H31").Value = 2321.3572
$ws.Range("I31").Value = 1249.9
$ws.Range("K31").Value = 1249.9
$ws.Range("M31").Value = -954.9000000000001
$ws.Range("H34").Value = 2321.3572
$ws.Range("I34").Value = 1249.9
$ws.Range("K34").Value = 1249.9
$ws.Range("M34").Value = -1047.9
$ws.Range("H113").Value = 2748.5
$ws.Range("I113").Value = 2631.5
$ws.Range("J113").Value = 3099.5
$ws.Range("K113").Value = 2631.5
$ws.Range("L113").Value = 3099.5
$ws.Range("M113").Value = -461.5
$ws.Range("N113").Value = -7439.5
$ws.Range("H122").Value = 2595.111
$ws.Range("I122").Value = 1904
$ws.Range("K122").Value = 5712
$ws.Range("M122").Value = -3262
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 3748
$ws.Range("I132").Value = 3748
$ws.Range("K132").Value = 11244
$ws.Range("M132").Value = -8714
$ws.Range("H134").Value = 1005.5
$ws.Range("I134").Value = 1005.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3016.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -481.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H101").Value = 7425
$ws.Range("J101").Value = 7552.6313
$ws.Range("L101").Value = 22657.8939
$ws.Range("N101").Value = -27525.8939

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H102").Value = 2200
$ws.Range("H122").Value = 4042.5715
$ws.Range("I122").Value = 4383
$ws.Range("K122").Value = 13149
$ws.Range("M122").Value = -10699
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 3841.3333
$ws.Range("I126").Value = 3841.3333
$ws.Range("K126").Value = 11523.9999
$ws.Range("M126").Value = -9053.999899999999
$ws.Range("H128").Value = 30774.166
$ws.Range("J128").Value = 30774.166
$ws.Range("L128").Value = 30774.166
$ws.Range("N128").Value = -40734.166
$ws.Range("H132").Value = 4270.6665
$ws.Range("I132").Value = 4998
$ws.Range("J132").Value = 3907
$ws.Range("K132").Value = 14994
$ws.Range("L132").Value = 11721
$ws.Range("M132").Value = -12464
$ws.Range("N132").Value = -16781

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1685.1428
$ws.Range("I7").Value = 1685.1428
$ws.Range("K7").Value = 1685.1428
$ws.Range("M7").Value = -1573.1428
$ws.Range("H40").Value = 3397.25
$ws.Range("I40").Value = 3397.25
$ws.Range("K40").Value = 3397.25
$ws.Range("M40").Value = -3261.25
$ws.Range("H122").Value = 6856.375
$ws.Range("I122").Value = 5169.2
$ws.Range("K122").Value = 15507.6
$ws.Range("M122").Value = -13057.6
$ws.Range("H126").Value = 1685.1428
$ws.Range("I126").Value = 1685.1428
$ws.Range("K126").Value = 5055.428400000001
$ws.Range("M126").Value = -2585.428400000001
$ws.Range("H132").Value = 3248.3333
$ws.Range("I132").Value = 1895
$ws.Range("J132").Value = 3925
$ws.Range("K132").Value = 5685
$ws.Range("L132").Value = 11775
$ws.Range("M132").Value = -3155
$ws.Range("N132").Value = -16835
$ws.Range("H136").Value = 957.8
$ws.Range("I136").Value = 947.25
$ws.Range("K136").Value = 2841.75
$ws.Range("M136").Value = -291.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45000
$ws.Range("J46").Value = 45000
$ws.Range("L46").Value = 45000
$ws.Range("N46").Value = -45462
$ws.Range("H122").Value = 3374
$ws.Range("I122").Value = 3374
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10122
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7672
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3667.9473
$ws.Range("I132").Value = 1919.9
$ws.Range("K132").Value = 5759.700000000001
$ws.Range("M132").Value = -3229.700000000001
$ws.Range("H134").Value = 45000
$ws.Range("J134").Value = 45000
$ws.Range("L134").Value = 135000
$ws.Range("N134").Value = -140070
$ws.Range("H136").Value = 1415.5333
$ws.Range("I136").Value = 1457.1538
$ws.Range("K136").Value = 4371.4614
$ws.Range("M136").Value = -1821.4614
